# added 4wk low sales check
# Updates the "Forecast Comparison" sheet (MyForecast flattens to 42 for the
# remaining forecast weeks, Inventory Coverage drains faster, Seasonality
# Index recalculated, and week 22 (row 13) now flagged "Urgent"), plus the
# corresponding roll-up totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# --- MyForecast (column D) ---------------------------------------------
$ws.Range("D3").Value  = 42
$ws.Range("D4").Value  = 42
$ws.Range("D5").Value  = 42
$ws.Range("D6").Value  = 42
$ws.Range("D7").Value  = 42
$ws.Range("D8").Value  = 42
$ws.Range("D9").Value  = 42
$ws.Range("D10").Value = 42
$ws.Range("D11").Value = 42
$ws.Range("D12").Value = 42
$ws.Range("D13").Value = 42
$ws.Range("D14").Value = 42
$ws.Range("D15").Value = 42
$ws.Range("D16").Value = 42
$ws.Range("D17").Value = 42

# --- Inventory Coverage (column H) --------------------------------------
$ws.Range("H3").Value  = 10.83
$ws.Range("H4").Value  = 9.83
$ws.Range("H5").Value  = 8.83
$ws.Range("H6").Value  = 7.83
$ws.Range("H7").Value  = 6.83
$ws.Range("H8").Value  = 5.83
$ws.Range("H9").Value  = 4.83
$ws.Range("H10").Value = 3.83
$ws.Range("H11").Value = 2.83
$ws.Range("H12").Value = 1.83
$ws.Range("H13").Value = 0.83
$ws.Range("H14").Value = 0

# --- Reorder Urgency (column J) -----------------------------------------
$ws.Range("J13").Value = "Urgent"

# --- Seasonality Index (column L) ---------------------------------------
$ws.Range("L2").Value  = 0.82
$ws.Range("L3").Value  = 1.09
$ws.Range("L4").Value  = 1.05
$ws.Range("L5").Value  = 0.98
$ws.Range("L6").Value  = 1.18
$ws.Range("L8").Value  = 1.07
$ws.Range("L9").Value  = 0.88
$ws.Range("L10").Value = 1.15
$ws.Range("L11").Value = 0.86
$ws.Range("L12").Value = 1.13
$ws.Range("L13").Value = 0.84
$ws.Range("L14").Value = 1.06
$ws.Range("L15").Value = 1.11
$ws.Range("L16").Value = 0.9
$ws.Range("L17").Value = 0.8

# --- Summary sheet roll-up totals ---------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").NumberFormat  = "@"
$ws2.Range("B9").Value  = "671"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "335"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "167"

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "41"
